$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @(0, 1, 1, 0, 2, 1, 0, 1, 1, 2, 2, 2, 0)

for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $i + 2
    $ws.Range("G$row").Value = $kValues[$i]
}
